$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5 and 6 swap places (the "Mac B. McGraw ..." record moves up to row 5,
# the "Neusha Barakati ..." record moves down to row 6), and the relocated
# row 5 record gets an updated affiliation string and cited_by_count.
#
# We move whole rows with Copy / PasteSpecial (values) rather than re-typing
# literal values, because several cells hold text that *looks* like a date
# or a plain integer (e.g. "2022-10-06", "2022", "0"); assigning such text
# directly to a cell lets Excel "smart" parse it into a real date/number,
# which would change its stored type away from the plain string it must stay.
# Copying an existing, already-correctly-typed cell's value is not re-parsed.

$xlPasteValues = -4163
$scratchRow = 50

# 1) Stash current row 5 ("Neusha Barakati ..." record) in a scratch row.
$ws.Range("A5:Q5").Copy() | Out-Null
$ws.Range("A" + $scratchRow + ":Q" + $scratchRow).PasteSpecial($xlPasteValues) | Out-Null

# 2) Move row 6 ("Mac B. McGraw ..." record) up into row 5.
$ws.Range("A6:Q6").Copy() | Out-Null
$ws.Range("A5:Q5").PasteSpecial($xlPasteValues) | Out-Null

# 3) Move the stashed original row 5 down into row 6.
$ws.Range("A" + $scratchRow + ":Q" + $scratchRow).Copy() | Out-Null
$ws.Range("A6:Q6").PasteSpecial($xlPasteValues) | Out-Null

# 4) Clean up the scratch row.
$ws.Range("A" + $scratchRow + ":Q" + $scratchRow).ClearContents() | Out-Null
$excel.CutCopyMode = 0

# 5) Update the relocated row 5 record's affiliation (plain text, safe to set directly).
$ws.Cells.Item(5, 2).Value2 = "The University of Arizona College of Medicine, United States; Exos, United States; The University of Arizona, United States; Arizona State University, United States; The University of Arizona, United States; The University of Arizona, United States"

# 6) Update the relocated row 5 record's cited_by_count to "1", forcing it to stay
#    a text value (like the rest of that column) instead of becoming a number.
$m5 = $ws.Cells.Item(5, 13)
$m5.NumberFormat = "@"
$m5.Value2 = "1"
$m5.ClearFormats()
